$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 133, pushing the existing rows 133:145 down to 134:146
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new weekly data point
$ws.Range("A133").Value = 1
$ws.Range("B133").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C133").Value = "Arica y Parinacota"
$ws.Range("D133").Value = 45218
$ws.Range("E133").Value = 15
$ws.Range("F133").Value = 100112038
$ws.Range("G133").Value = "Cebollín baby"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 650
$ws.Range("K133").Value = 1400
$ws.Range("L133").Value = 1500
$ws.Range("M133").Value = 1454
$ws.Range("N133").Value = '$/paquete 1,5 a 2 kilos'
$ws.Range("O133").Value = "Región de Arica y Parinacota"
$ws.Range("P133").Value = 727
$ws.Range("Q133").Value = 2
$ws.Range("R133").Value = "Hortaliza"
